$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows (Costa Rica, month 2022-02) to be appended below the existing
# 2022-01 data block, mirroring the same Subsector order.
$data = @(
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Education",88,0,88,88,0,0,0,0,2,0,68,18,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Food Security",47,41,56,56,0,0,0,0,8,9,20,19,0,0,50,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Health",1264,0,183,183,0,0,0,0,8,5,117,53,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Humanitarian Transportation",0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Integration",79,0,115,115,0,0,0,0,3,0,83,29,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Multipurpose Cash Assistance (MPC)",83,82,84,84,0,0,0,0,18,12,33,21,0,0,83,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Nutrition",0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Protection (General)",172,0,296,296,0,0,0,0,7,5,175,108,0,1,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Protection (GBV)",1,0,3,3,0,0,0,0,0,0,3,0,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Protection (Human Trafficking and Smuggling)",0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Protection (Child Protection)",0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Shelter",59,0,60,60,0,0,0,0,8,8,15,29,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","WASH",42,0,52,52,0,0,0,0,9,8,18,17,0,0,0,"Ok","Ok","ok"),
        @("Central America and Mexico","Costa Rica","Country level","2022-02","Intersector",1835,123,937,937,0,0,0,0,63,47,532,294,0,1,133,"Ok","Ok","ok")
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}
